# Wijzigingen en aanvullingen prognoses.
# - Rename the worksheet from "prog002" to "data"
# - Configure the page setup (paper size A4 = 9, portrait orientation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "data"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
